$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2, shifting the existing data rows (old 2..14) down to (3..15)
$ws.Rows.Item(2).Insert()

# Excel copies the formatting of the row above (the header) into the newly inserted row;
# clear that so the new row starts out with the default (unstyled) formatting like the
# other data rows.
$ws.Rows.Item(2).ClearFormats()

# Populate the newly inserted row 2 with the new weekly data record
$ws.Cells.Item(2, 1).Value = 4
$ws.Cells.Item(2, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(2, 3).Value = "Los Lagos"
$ws.Cells.Item(2, 4).Value = 44547
$ws.Cells.Item(2, 5).Value = 10
$ws.Cells.Item(2, 6).Value = "Fruta"
$ws.Cells.Item(2, 7).Value = 100101
$ws.Cells.Item(2, 8).Value = "Berries"
$ws.Cells.Item(2, 9).Value = 100101001
$ws.Cells.Item(2, 10).Value = "Arándano (blue)"
$ws.Cells.Item(2, 11).Value = "Sin especificar"
$ws.Cells.Item(2, 12).Value = "Primera"
$ws.Cells.Item(2, 13).Value = 400
$ws.Cells.Item(2, 14).Value = 5000
$ws.Cells.Item(2, 15).Value = 5500
$ws.Cells.Item(2, 16).Value = 5250
$ws.Cells.Item(2, 17).Value = "$/bandeja 12 canastillos 125 gramos"
$ws.Cells.Item(2, 18).Value = "Región del Maule"
$ws.Cells.Item(2, 19).Value = 3500
$ws.Cells.Item(2, 20).Value = 1.5

# Re-apply the date number format to the Fecha column (D), matching the other rows,
# by copying the style from the row below (which holds the original row-2 formatting).
$ws.Cells.Item(3, 4).Copy()
$ws.Cells.Item(2, 4).PasteSpecial(-4122)
$ws.Cells.Item(2, 4).Value = 44547
